# Add a new "free prog users" entry (row 8) with an email hyperlink,
# mirroring the existing "free prog users" row (row 7 / B1 hyperlink style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label / value for row 8
$ws.Range("A8").Value = "add new users"
$ws.Range("B8").Value = "rishabh.singh+3@snackmagic.com"

# Turn B8 into a mailto: hyperlink (adds a new relationship + hyperlink entry)
[void]$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:rishabh.singh+3@snackmagic.com")

# Match the existing hyperlink cell style (wrapped text + Hyperlink style),
# same as B1. Order matters: apply Style first, then WrapText, so the engine
# reuses the existing style slot instead of allocating a new one.
$ws.Range("B8").Style = "Hyperlink"
$ws.Range("B8").WrapText = $true

# Update selection to the newly added cell
[void]$ws.Range("B8").Select()
